# Generate Report for handback
# This script updates the zh-cn and de-de handback status sheets:
#  - updates the "Status" column text to reflect a completed handback
#  - fills in the "Latest Target File" / "Latest Handback File" columns (which were
#    previously blank) with hyperlinked file names
#  - stamps the "Latest Handback DateTime" column with the handback timestamp

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

function Update-LangSheet {
    param(
        [string]$SheetName,
        [string]$XlfFileName,
        [string]$HandbackDateTime
    )

    $ws = $wb.Worksheets.Item($SheetName)

    $aUrl = "https://github.com/OpenLocalizationTest/oltest/blob/c23857ddba523431bcd1eec2e66c903c973b86c7/e2e/a.md.md"
    $xlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fb502e11503e052d43f06d9689c1b9f2d27585e3/ol-handoff/OpenLocalizationTestOrg/oltest.$SheetName/yuwzho/$XlfFileName"
    if ($SheetName -eq "de-de") {
        $xlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7909717437d8a48f542c0f3b937272bd596902d3/ol-handoff/OpenLocalizationTestOrg/oltest.$SheetName/yuwzho/$XlfFileName"
    }

    foreach ($row in 2, 3) {
        # Status column (B) -> handed back
        $ws.Range("B$row").Value = $newStatus

        # Latest Target File (E) - same source file as column A
        $ws.Range("E$row").Value = "a.md.md"
        $ws.Hyperlinks.Add($ws.Range("E$row"), $aUrl, "", "", "a.md.md")

        # Latest Handback File (F) - the handed-back xlf, same as column C
        $ws.Range("F$row").Value = $XlfFileName
        $ws.Hyperlinks.Add($ws.Range("F$row"), $xlfUrl, "", "", $XlfFileName)

        # Latest Handback DateTime (G)
        $ws.Range("G$row").Value = $HandbackDateTime
    }
}

Update-LangSheet "zh-cn" "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf" "2016-01-26 05:34:10"
Update-LangSheet "de-de" "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf" "2016-01-26 05:34:27"
